# Fruta / hortaliza, semanal
# Insert one new weekly data row at row 37 (right after the header block of
# fixed rows 2-36), pushing the existing rows 37-237 down to 38-238.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 37; Excel shifts rows 37..237 down to 38..238
# and extends the used range (dimension becomes A1:R238 automatically).
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with the new weekly observation.
$ws.Range("A37").Value = 3
$ws.Range("B37").Value = "Femacal de La Calera"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44561
$ws.Range("E37").Value = 5
$ws.Range("F37").Value = 100112039
$ws.Range("G37").Value = "Ciboulette"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 180
$ws.Range("K37").Value = 1500
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = 1500
$ws.Range("N37").Value = '$/docena de atados'
$ws.Range("O37").Value = "Provincia de Quillota"
$ws.Range("P37").Value = 500
$ws.Range("Q37").Value = 3
$ws.Range("R37").Value = "Hortaliza"
